$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that needs to move
# from 45181 (2023-09-12) to 45182 (2023-09-13) for every data row (2-135).
for ($row = 2; $row -le 135; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
